$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 240.97
$ws.Range("C2").Value = 192.78
$ws.Range("D2").Value = 289.17

$ws.Range("B3").Value = 360.54
$ws.Range("C3").Value = 288.43
$ws.Range("D3").Value = 432.64

$ws.Range("B4").Value = 274.69
$ws.Range("C4").Value = 219.75
$ws.Range("D4").Value = 329.62
